$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 230 (2019-11-15) with revised figures ---
$ws.Range("B230").Value = 5827
$ws.Range("C230").Value = 5709
$ws.Range("D230").Value = 3568
$ws.Range("E230").Value = 1768

# --- Fill in row 231 (2019-12-15), previously blank placeholders ---
$ws.Range("B231").Value = 5907
$ws.Range("C231").Value = 5730
$ws.Range("D231").Value = 3488
$ws.Range("E231").Value = 1895

# --- Row 232 (2020-01-15): new "NA" placeholders in B:E ---
$ws.Range("B232:E232").Value = "NA"

# --- Row 233 (2020-02-15): brand-new row ---
# Copy the date cell above so the new date cell inherits the same date
# number format / style, then overwrite with the new date value.
$ws.Range("A232").Copy($ws.Range("A233"))
$ws.Range("A233").Value = 43876
$ws.Range("B233:E233").Value = "NA"

# --- Update the "TRAFO" data-range bookmark stored on A1's validation ---
$dv = $ws.Range("A1").Validation
$dv.InputMessage = '$A$1:$E$233'
